$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6, shifting existing rows (6-31) down to (7-32)
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with a fresh weekly record (copy of row 7's
# category/label fields, new date + updated price figures)
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(6, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(6, 4).Value = [DateTime]"2022-05-06"
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(6, 6).Value = "Fruta"
$ws.Cells.Item(6, 7).Value = 100108
$ws.Cells.Item(6, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(6, 9).Value = 100108001
$ws.Cells.Item(6, 10).Value = "Guayaba"
$ws.Cells.Item(6, 11).Value = "Sin especificar"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 120
$ws.Cells.Item(6, 14).Value = 1300
$ws.Cells.Item(6, 15).Value = 1400
$ws.Cells.Item(6, 16).Value = 1350
$ws.Cells.Item(6, 17).Value = "$/kilo (en caja de 10 kilos )"
$ws.Cells.Item(6, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(6, 19).Value = 1350
$ws.Cells.Item(6, 20).Value = 1
